## CISC 3140 Lab 6 - "added log of repo 2"
##
## The sentence:
##   "Also cause I cant seem to get the log to not be unreadable at times: I "
## becomes:
##   "Also because I can't seem to get the log to not be unreadable at times: I "
## (word "be" inserted before "cause", curly apostrophe inserted in "can't")
## and - per the target OOXML - the single run that used to hold this whole
## sentence is split into five separate runs along the edited word
## boundaries, exactly as Word would leave it after an in-place typing edit.

$d = $word.ActiveDocument

$oldText = "Also cause I cant seem to get the log to not be unreadable at times: I "
$curlyApostrophe = [char]0x2019
$newText = "Also because I can" + $curlyApostrophe + "t seem to get the log to not be unreadable at times: I "

# --- Step 1: fix the wording/spelling in place -----------------------------
$fixRange = $d.Content
$found = $fixRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) {
    throw "Could not locate the target sentence to edit"
}

# --- Step 2: re-locate the corrected sentence so we know its start offset --
$locateRange = $d.Content
$locateRange.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentenceStart = $locateRange.Start

# --- Step 3: split the single run into the five runs the diff expects ------
# "Also " | "be" | "cause I can" | "'" | "t seem to get the log to not be unreadable at times: I "
# Dropping a zero-width bookmark at each internal boundary (and immediately
# deleting it) forces the engine to keep the runs distinct instead of
# re-coalescing them because they share identical (empty) run formatting.
$segmentLengths = @(5, 2, 11, 1)

$bookmarkNames = New-Object System.Collections.ArrayList
$cursor = $sentenceStart
$idx = 0
foreach ($segLen in $segmentLengths) {
    $cursor = $cursor + $segLen
    $idx = $idx + 1
    $bmName = "tmp_split_$idx"
    $pt = $d.Range($cursor, $cursor)
    $d.Bookmarks.Add($bmName, $pt) | Out-Null
    $bookmarkNames.Add($bmName) | Out-Null
}

# Also drop one right after the edited sentence so the trailing edited run
# ("t seem ... I ") stays separate from the untouched "copy/pasted..." run
# that follows it.
$idx = $idx + 1
$bmName = "tmp_split_$idx"
$tailPoint = $d.Range($sentenceStart + $newText.Length, $sentenceStart + $newText.Length)
$d.Bookmarks.Add($bmName, $tailPoint) | Out-Null
$bookmarkNames.Add($bmName) | Out-Null

foreach ($bmName in $bookmarkNames) {
    $d.Bookmarks($bmName).Delete()
}

Write-Output "Updated sentence now reads:"
Write-Output $locateRange.Text
